$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 89

$ws.Range("A$row").Value = 44824
$ws.Range("A$row").NumberFormat = "YYYY-MM-DD"

$ws.Range("B$row").Value = "LICI"
$ws.Range("C$row").Value = "EQ"

$ws.Range("D$row").Value = 654.75
$ws.Range("E$row").Value = 657.55
$ws.Range("F$row").Value = 665.05
$ws.Range("G$row").Value = 655.55
$ws.Range("H$row").Value = 656.4
$ws.Range("I$row").Value = 656.35
$ws.Range("J$row").Value = 659.3200000000001
$ws.Range("K$row").Value = 927521
$ws.Range("L$row").Value = 61153015895000.01
$ws.Range("M$row").Value = 28388
$ws.Range("N$row").Value = 393288
$ws.Range("O$row").Value = 0.424
